$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "298.26"); assigning
# them directly via .Value would make Excel auto-coerce them to real numbers.
# Forcing a Text number format keeps them as strings like the source data, and
# resetting the style back to Normal afterwards avoids leaving a stray style index
# on the cell (matching cells that never had an explicit style).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "42.262.78"
$ws.Range("E2").Value = "  -1.76%  "

Set-TextValue $ws.Range("D3") "2.271.55"
$ws.Range("E3").Value = "  -2.84%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue $ws.Range("D5") "298.26"
$ws.Range("E5").Value = "  -2.76%  "

Set-TextValue $ws.Range("D6") "95.43"
$ws.Range("E6").Value = "  -5.64%  "

$ws.Range("E7").Value = "  +0.05%  "

Set-TextValue $ws.Range("D8") "0.494"
$ws.Range("E8").Value = "  -3.72%  "

$ws.Range("E9").Value = "  -3.79%  "

Set-TextValue $ws.Range("D10") "33.42"
$ws.Range("E10").Value = "  -4.29%  "

$ws.Range("E11").Value = "  -1.03%  "

Set-TextValue $ws.Range("D12") "48.68"
$ws.Range("E12").Value = "  -7.31%  "

$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("E14").Value = "  -3.00%  "

Set-TextValue $ws.Range("D15") "15.72"
$ws.Range("E15").Value = "  -0.58%  "

Set-TextValue $ws.Range("D16") "2.625.90"
$ws.Range("E16").Value = "  -2.95%  "

Set-TextValue $ws.Range("D17") "2.272.92"
$ws.Range("E17").Value = "  -3.02%  "

Set-TextValue $ws.Range("D18") "0.780"
$ws.Range("E18").Value = "  -5.89%  "

Set-TextValue $ws.Range("D19") "42.198.49"
$ws.Range("E19").Value = "  -1.75%  "

Set-TextValue $ws.Range("D20") "11.64"
$ws.Range("E20").Value = "  -0.77%  "

Set-TextValue $ws.Range("D21") "0.0₃0890"
$ws.Range("E21").Value = "  -2.15%  "

Set-TextValue $ws.Range("D22") "6.00"
$ws.Range("E22").Value = "  -3.64%  "

Set-TextValue $ws.Range("D23") "66.73"
$ws.Range("E23").Value = "  -1.98%  "

Set-TextValue $ws.Range("D24") "233.42"
$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("E25").Value = "  -2.54%  "

$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("E27").Value = "  -4.00%  "

Set-TextValue $ws.Range("D28") "23.89"
$ws.Range("E28").Value = "  -7.00%  "

$ws.Range("E29").Value = "  -1.14%  "

Set-TextValue $ws.Range("D30") "167.59"
$ws.Range("E30").Value = "  +3.02%  "

Set-TextValue $ws.Range("D31") "34.08"
$ws.Range("E31").Value = "  -4.96%  "

Set-TextValue $ws.Range("D32") "9.11"
$ws.Range("E32").Value = "  -2.26%  "

$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("E34").Value = "  -4.01%  "

$ws.Range("E35").Value = "  -2.31%  "

Set-TextValue $ws.Range("D36") "0.0691"
$ws.Range("E36").Value = "  -4.91%  "

$ws.Range("E37").Value = "  -4.91%  "

Set-TextValue $ws.Range("D38") "16.41"
$ws.Range("E38").Value = "  -6.59%  "

$ws.Range("E39").Value = "  -4.09%  "

$ws.Range("E40").Value = "  -3.01%  "

$ws.Range("E41").Value = "  -3.22%  "

$ws.Range("E42").Value = "  -6.75%  "

$ws.Range("E43").Value = "  -7.48%  "

Set-TextValue $ws.Range("D44") "1.966.03"
$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("E45").Value = "  -2.20%  "

Set-TextValue $ws.Range("D46") "17.56"
$ws.Range("E46").Value = "  -7.27%  "

Set-TextValue $ws.Range("D47") "9.60"
$ws.Range("E47").Value = "  -5.51%  "

$ws.Range("E48").Value = "  -4.95%  "

Set-TextValue $ws.Range("D49") "2.498.09"
$ws.Range("E49").Value = "  -2.37%  "

Set-TextValue $ws.Range("D50") "52.28"
$ws.Range("E50").Value = "  -7.88%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D51") "4.57"
$ws.Range("E51").Value = "  -2.95%  "
